$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - Division U7 (Tier: 1), Post-Optimization
$ws.Range("C16").Value = 38.0
$ws.Range("D16").Value = 6.0
$ws.Range("F16").Value = 8.0
$ws.Range("G16").Replace("79.55%", "81.82%")
$ws.Range("H16").Replace("84.09%", "86.36%")
$ws.Range("K16").Value = 92.55

# Row 17 - Division U7 (Tier: 2), Post-Optimization
$ws.Range("K17").Value = 120.18

# Row 18 - Division U7 (Tier: 3), Post-Optimization
$ws.Range("C18").Value = 78.0
$ws.Range("D18").Value = 10.0
$ws.Range("H18").Replace("87.5%", "88.64%")
$ws.Range("K18").Value = 101.91

# Row 20 - Division U8 (Tier: 1), Post-Optimization
$ws.Range("F20").Value = 11.0
$ws.Range("G20").Replace("45.45%", "50.0%")
$ws.Range("K20").Value = 157.0

# Row 21 - Division U8 (Tier: 2), Post-Optimization
$ws.Range("F21").Value = 10.0
$ws.Range("G21").Replace("81.82%", "77.27%")
$ws.Range("K21").Value = 93.18

# Row 23 - Division U9 (Tier: 0), Post-Optimization
$ws.Range("K23").Value = 107.82

# Row 24 - Division U9 (Tier: 1), Post-Optimization
$ws.Range("F24").Value = 11.0
$ws.Range("G24").Replace("72.73%", "75.0%")
$ws.Range("K24").Value = 100.0

# Row 25 - Division U9 (Tier: 2), Post-Optimization
$ws.Range("C25").Value = 36.0
$ws.Range("D25").Value = 8.0
$ws.Range("F25").Value = 9.0
$ws.Range("G25").Replace("75.0%", "79.55%")
$ws.Range("H25").Replace("84.09%", "81.82%")
$ws.Range("K25").Value = 113.36

# Row 26 - Division U9 (Tier: 3), Post-Optimization
$ws.Range("C26").Value = 114.0
$ws.Range("D26").Value = 18.0
$ws.Range("F26").Value = 41.0
$ws.Range("G26").Replace("69.7%", "68.94%")
$ws.Range("H26").Replace("84.85%", "86.36%")
$ws.Range("K26").Value = 209.18

# Row 30 - Optimization Time value
$ws.Range("A30").Replace("3 min, 22 sec", "2 min, 59 sec")
